$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextCell $ws.Range("D2") "61.354.07"
Set-TextCell $ws.Range("E2") "  -1.91%  "

Set-TextCell $ws.Range("D3") "3.391.71"
Set-TextCell $ws.Range("E3") "  +0.81%  "

Set-TextCell $ws.Range("E4") "  -0.11%  "

Set-TextCell $ws.Range("D5") "574.60"
Set-TextCell $ws.Range("E5") "  +1.25%  "

Set-TextCell $ws.Range("D6") "136.43"
Set-TextCell $ws.Range("E6") "  +9.60%  "

Set-TextCell $ws.Range("E7") "  -0.06%  "

Set-TextCell $ws.Range("D8") "3.391.17"
Set-TextCell $ws.Range("E8") "  +0.89%  "

Set-TextCell $ws.Range("D9") "0.478"
Set-TextCell $ws.Range("E9") "  +1.61%  "

Set-TextCell $ws.Range("D10") "7.59"
Set-TextCell $ws.Range("E10") "  +3.92%  "

Set-TextCell $ws.Range("E11") "  +4.40%  "

Set-TextCell $ws.Range("D12") "0.387"
Set-TextCell $ws.Range("E12") "  +4.00%  "

Set-TextCell $ws.Range("D13") "3.970.28"
Set-TextCell $ws.Range("E13") "  +0.21%  "

Set-TextCell $ws.Range("E14") "  +1.08%  "

Set-TextCell $ws.Range("E15") "  +3.03%  "

Set-TextCell $ws.Range("D16") "3.389.05"
Set-TextCell $ws.Range("E16") "  +0.15%  "

Set-TextCell $ws.Range("D17") "25.20"
Set-TextCell $ws.Range("E17") "  +3.69%  "

Set-TextCell $ws.Range("D18") "61.494.00"
Set-TextCell $ws.Range("E18") "  -1.71%  "

Set-TextCell $ws.Range("D19") "14.10"
Set-TextCell $ws.Range("E19") "  +8.64%  "

Set-TextCell $ws.Range("E20") "  +1.00%  "

Set-TextCell $ws.Range("D21") "5.79"
Set-TextCell $ws.Range("E21") "  +3.59%  "

Set-TextCell $ws.Range("D22") "377.32"
Set-TextCell $ws.Range("E22") "  +1.81%  "

Set-TextCell $ws.Range("E23") "  +3.55%  "

Set-TextCell $ws.Range("D24") "3.525.26"
Set-TextCell $ws.Range("E24") "  +0.68%  "

Set-TextCell $ws.Range("E25") "  +0.25%  "

Set-TextCell $ws.Range("D26") "70.78"
Set-TextCell $ws.Range("E26") "  -0.81%  "

Set-TextCell $ws.Range("E27") "  +13.01%  "

Set-TextCell $ws.Range("E28") "  +23.92%  "

Set-TextCell $ws.Range("E29") "  +13.67%  "

Set-TextCell $ws.Range("E30") "  +0.51%  "

Set-TextCell $ws.Range("D31") "8.17"
Set-TextCell $ws.Range("E31") "  +6.77%  "

Set-TextCell $ws.Range("E32") "  +2.08%  "

Set-TextCell $ws.Range("E33") "  +5.01%  "

Set-TextCell $ws.Range("E34") "  -0.05%  "

Set-TextCell $ws.Range("D35") "3.423.85"
Set-TextCell $ws.Range("E35") "  +0.90%  "

Set-TextCell $ws.Range("D36") "23.46"
Set-TextCell $ws.Range("E36") "  +3.93%  "

Set-TextCell $ws.Range("E37") "  +8.74%  "

Set-TextCell $ws.Range("D38") "1.58"
Set-TextCell $ws.Range("E38") "  +8.16%  "

Set-TextCell $ws.Range("D39") "6.95"
Set-TextCell $ws.Range("E39") "  +5.84%  "

Set-TextCell $ws.Range("D40") "163.06"
Set-TextCell $ws.Range("E40") "  -0.94%  "

Set-TextCell $ws.Range("D41") "0.0790"
Set-TextCell $ws.Range("E41") "  +6.28%  "

Set-TextCell $ws.Range("E42") "  -0.26%  "

Set-TextCell $ws.Range("E43") "  +15.48%  "

Set-TextCell $ws.Range("D44") "4.42"
Set-TextCell $ws.Range("E44") "  +5.45%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D45") "41.49"
Set-TextCell $ws.Range("E45") "  +0.56%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D46") "0.762"
Set-TextCell $ws.Range("E46") "  -0.75%  "

Set-TextCell $ws.Range("E47") "  +5.85%  "

Set-TextCell $ws.Range("D48") "23.49"
Set-TextCell $ws.Range("E48") "  +4.13%  "

Set-TextCell $ws.Range("E49") "  +6.26%  "

Set-TextCell $ws.Range("D50") "23.11"
Set-TextCell $ws.Range("E50") "  +16.01%  "

Set-TextCell $ws.Range("D51") "0.902"
Set-TextCell $ws.Range("E51") "  +7.03%  "
